$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("G7").Value = 9.75
$ws.Range("H7").Value = 6.5
$ws.Range("I7").Value = 3.25

$ws.Range("H10").Select()
